# Generate Report for Archive
#
# Two source files (73603af9-183d-4d23-8044-70b790fdaa8e.md and
# c44e843b-c8b2-4011-a08e-70b02a849b95.md) have moved from "Ready for
# handoff" to "In Translation" in the localization status report.
# Update the Status columns on the Overview sheet (rows 8 and 9, columns
# B and C) as well as on the per-locale "zh-cn" and "de-de" sheets (rows
# 8 and 9, column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C8").Value = "In Translation"
$zhcn.Range("C9").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C8").Value = "In Translation"
$dede.Range("C9").Value = "In Translation"
